$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect to update the cells, then
# re-protect it afterwards so the workbook ends up protected again.
$ws.Unprotect("D382")

# Update the confidential disclaimer text: date 2021-05-03 -> 2021-05-04
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.09961281441069632
$ws.Range("E2").Value = -0.03549543299803448

$ws.Range("D3").Value = 0.1082001511691157
$ws.Range("E3").Value = -0.01952873764836582

$ws.Range("D4").Value = 0.1185739615328938
$ws.Range("E4").Value = -0.004522769806612703

$ws.Range("D5").Value = 0.1392174402415049
$ws.Range("E5").Value = -0.004753981459472256

$ws.Range("D6").Value = 0.1350415399397936
$ws.Range("E6").Value = -0.004006078187595019

$ws.Range("D7").Value = 0.1432507535733523
$ws.Range("E7").Value = -0.00171216588985057

$ws.Range("D8").Value = 0.1271132826564103
$ws.Range("E8").Value = -0.01761080129145887

$ws.Range("D9").Value = 0.1289900564762331
$ws.Range("E9").Value = -0.01528254614457125

$ws.Range("D10").Value = 1
$ws.Range("E10").Value = -0.0118430514823622

# Restore sheet protection (matches original: objects/scenarios protected,
# column/row formatting allowed).
$ws.Protect("D382", $true, $true, $true, $false, $false, $true, $true)
